$wb = $excel.ActiveWorkbook

# --- 1. Rename "Requested quantity" header on existing sheets ---
$weekly = $wb.Sheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Sheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$forecast.Name = "PO Forecast"

# Copy the header formatting (bold font, border, center/top alignment) from
# the Weekly Quantity header row so the new header cells share the same style.
$weekly.Range("A1:B1").Copy()
$forecast.Range("A1:D1").PasteSpecial(-4122)

$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

# Copy the date-cell formatting (YYYY-MM-DD HH:MM:SS number format) from the
# existing sheet so the new "ds" column matches it.
$weekly.Range("A2").Copy()
$forecast.Range("A2:A16").PasteSpecial(-4122)

$rows = @(
    @(45011.99999999999, 31, 6.229027425561182, 54.56085802102525),
    @(45025.99999999999, 39, 14.62629175573134, 61.62876921686846),
    @(45039.99999999999, 46, 22.96135995208966, 69.25778382902132),
    @(45046.99999999999, 50, 27.10312485182257, 74.83711115352112),
    @(45053.99999999999, 54, 29.61953544965583, 75.72322403684733),
    @(45060.99999999999, 58, 35.45673321370265, 79.83588084553408),
    @(45067.99999999999, 62, 39.14395884808096, 83.43106226309268),
    @(45074.99999999999, 66, 40.88012843941631, 88.39186979994645),
    @(45081.99999999999, 70, 45.93134808352851, 92.78398823224195),
    @(45088.99999999999, 74, 51.40862558445972, 97.44021547820476),
    @(45095.99999999999, 77, 54.69106259185307, 99.88120448903834),
    @(45102.99999999999, 81, 58.12508903754218, 104.4928196636394),
    @(45109.99999999999, 85, 60.45573279901129, 107.0855195643106),
    @(45116.99999999999, 89, 65.40880373071123, 112.3466885313853),
    @(45123.99999999999, 93, 69.3807902779148, 116.2028208023648)
)

$r = 2
foreach ($row in $rows) {
    $forecast.Cells.Item($r, 1).Value = $row[0]
    $forecast.Cells.Item($r, 2).Value = $row[1]
    $forecast.Cells.Item($r, 3).Value = $row[2]
    $forecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
